$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "[81, 30, 42, 49, 25]"
$ws.Range("B7").Value = 0.9801307842973825
$ws.Range("B8").Value = 0.0003852696616908347
$ws.Range("B9").Value = 62
$ws.Range("B10").Value = 71
$ws.Range("B12").Value = "[[30, 1], [84, 41]]"
$ws.Range("B14").Value = "[[57, 78, 22], [21, 91, 64], [22, 13, 18]]"
$ws.Range("B15").Value = "[0.9863555691481738, 0.6768682690872199, 0.6892312689820487, 0.882524473202599]"
